$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 34487
$ws.Range("I40").Value = 41666
$ws.Range("K40").Value = 41666
$ws.Range("M40").Value = -41491

$ws.Range("H96").Value = 1187.6875
$ws.Range("I96").Value = 372.42856
$ws.Range("J96").Value = 1821.7778
$ws.Range("K96").Value = 1117.28568
$ws.Range("L96").Value = 5465.3334
$ws.Range("M96").Value = 255.71432
$ws.Range("N96").Value = -8211.3334

$ws.Range("H132").Value = 104509.17
$ws.Range("I132").Value = 242130.73
$ws.Range("J132").Value = 13580.643
$ws.Range("K132").Value = 726392.1900000001
$ws.Range("L132").Value = 40741.929
$ws.Range("M132").Value = -723862.1900000001
$ws.Range("N132").Value = -45801.929

$ws.Range("H135").Value = 2271.6
$ws.Range("I135").Value = 698
$ws.Range("K135").Value = 6282
$ws.Range("M135").Value = -3747

$ws.Range("H137").Value = 4131.5835
$ws.Range("I137").Value = 1732.3334
$ws.Range("J137").Value = 4931.3335
$ws.Range("K137").Value = 5197.0002
$ws.Range("L137").Value = 14794.0005
$ws.Range("M137").Value = -2647.0002
$ws.Range("N137").Value = -19894.0005

$ws.Range("H138").Value = 7943.212
$ws.Range("J138").Value = 8241.161
$ws.Range("L138").Value = 24723.483
$ws.Range("N138").Value = -35003.483

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2010.0667
$ws.Range("I45").Value = 1994.9286
$ws.Range("K45").Value = 1994.9286
$ws.Range("M45").Value = -1617.9286

$ws.Range("H74").Value = 1473.4706
$ws.Range("I74").Value = 774.4
$ws.Range("K74").Value = 774.4
$ws.Range("M74").Value = 99.60000000000002

$ws.Range("H77").Value = 1473.4706
$ws.Range("I77").Value = 774.4
$ws.Range("K77").Value = 3872
$ws.Range("M77").Value = 496

$ws.Range("H102").Value = 762157.25
$ws.Range("I102").Value = 1371062
$ws.Range("K102").Value = 1371062
$ws.Range("M102").Value = -1369440

$ws.Range("H104").Value = 40166.668
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()

$ws.Range("H106").Value = 83000
$ws.Range("J106").Value = 83000
$ws.Range("L106").Value = 83000
$ws.Range("N106").Value = -85524

$ws.Range("H122").Value = 3147.037
$ws.Range("I122").Value = 2246.7222
$ws.Range("K122").Value = 6740.1666
$ws.Range("M122").Value = -4290.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 336.25
$ws.Range("J80").Value = 364
$ws.Range("L80").Value = 364
$ws.Range("N80").Value = -2360

$ws.Range("H83").Value = 336.25
$ws.Range("J83").Value = 364
$ws.Range("L83").Value = 1820
$ws.Range("N83").Value = -11804

$ws.Range("H107").Value = 1735.238
$ws.Range("I107").Value = 2639.625
$ws.Range("K107").Value = 2639.625
$ws.Range("M107").Value = -719.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2733.842
$ws.Range("I31").Value = 920.5
$ws.Range("J31").Value = 5842.4287
$ws.Range("K31").Value = 920.5
$ws.Range("L31").Value = 5842.4287
$ws.Range("M31").Value = -625.5
$ws.Range("N31").Value = -6432.4287

$ws.Range("H34").Value = 2733.842
$ws.Range("I34").Value = 920.5
$ws.Range("J34").Value = 5842.4287
$ws.Range("K34").Value = 920.5
$ws.Range("L34").Value = 5842.4287
$ws.Range("M34").Value = -718.5
$ws.Range("N34").Value = -6246.4287

$ws.Range("H64").Value = 18999.5
$ws.Range("J64").Value = 18999.5
$ws.Range("L64").Value = 18999.5
$ws.Range("N64").Value = -19495.5

$ws.Range("H67").Value = 18999.5
$ws.Range("J67").Value = 18999.5
$ws.Range("L67").Value = 18999.5
$ws.Range("N67").Value = -20715.5

$ws.Range("H107").Value = 1013907.75
$ws.Range("I107").Value = 1655063.5
$ws.Range("J107").Value = 6377.143
$ws.Range("K107").Value = 1655063.5
$ws.Range("L107").Value = 6377.143
$ws.Range("M107").Value = -1653143.5
$ws.Range("N107").Value = -10217.143

$ws.Range("H132").Value = 7250
$ws.Range("I132").Value = 3540
$ws.Range("J132").Value = 38166.668
$ws.Range("K132").Value = 10620
$ws.Range("L132").Value = 114500.004
$ws.Range("M132").Value = -8090
$ws.Range("N132").Value = -119560.004

$ws.Range("H134").Value = 2118.9148
$ws.Range("I134").Value = 2147.3809
$ws.Range("K134").Value = 6442.1427
$ws.Range("M134").Value = -3907.1427

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1833.1666
$ws.Range("I98").Value = 1499.5
$ws.Range("K98").Value = 4498.5
$ws.Range("M98").Value = -3000.5

$ws.Range("H128").Value = 349994
$ws.Range("I128").Value = 349994
$ws.Range("K128").Value = 1049982
$ws.Range("M128").Value = -1045002

$ws.Range("H131").Value = 8066946.5
$ws.Range("I131").Value = 12821260
$ws.Range("K131").Value = 38463780
$ws.Range("M131").Value = -38458740

$ws.Range("H133").Value = 17387.082
$ws.Range("I133").Value = 24267
$ws.Range("K133").Value = 72801
$ws.Range("M133").Value = -67741

$ws.Range("H134").Value = 37454.332
$ws.Range("I134").Value = 37454.332
$ws.Range("K134").Value = 112362.996
$ws.Range("M134").Value = -107292.996

$ws.Range("H138").Value = 66376.5
$ws.Range("I138").Value = 146685.72
$ws.Range("J138").Value = 3913.7778
$ws.Range("K138").Value = 440057.16
$ws.Range("L138").Value = 11741.3334
$ws.Range("M138").Value = -434917.16
$ws.Range("N138").Value = -22021.3334

$ws.Range("H139").Value = 837377.2
$ws.Range("I139").Value = 1113192
$ws.Range("K139").Value = 3339576
$ws.Range("M139").Value = -3334436

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1678599.8
$ws.Range("I80").Value = 5544663
$ws.Range("K80").Value = 5544663
$ws.Range("M80").Value = -5543665

$ws.Range("H83").Value = 1678599.8
$ws.Range("I83").Value = 5544663
$ws.Range("K83").Value = 27723315
$ws.Range("M83").Value = -27718323

$ws.Range("H113").Value = 9984.875
$ws.Range("I113").Value = 6119.75
$ws.Range("K113").Value = 6119.75
$ws.Range("M113").Value = -3949.75

$ws.Range("H122").Value = 5501203
$ws.Range("I122").Value = 5501203
$ws.Range("K122").Value = 16503609
$ws.Range("M122").Value = -16501159

$ws.Range("H132").Value = 2662
$ws.Range("I132").Value = 2662
$ws.Range("K132").Value = 7986
$ws.Range("M132").Value = -5456

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 32500
$ws.Range("I2").Value = 32500
$ws.Range("K2").Value = 32500
$ws.Range("M2").Value = -32388

$ws.Range("H61").Value = 4951.4287
$ws.Range("I61").Value = 5226.6665
$ws.Range("K61").Value = 5226.6665
$ws.Range("M61").Value = -5024.6665

$ws.Range("H76").Value = 41253.285
$ws.Range("J76").Value = 41253.285
$ws.Range("L76").Value = 41253.285
$ws.Range("N76").Value = -41929.285

$ws.Range("H79").Value = 41253.285
$ws.Range("J79").Value = 41253.285
$ws.Range("L79").Value = 41253.285
$ws.Range("N79").Value = -43593.285

$ws.Range("H82").Value = 2233439.8
$ws.Range("J82").Value = 1781.6666
$ws.Range("L82").Value = 1781.6666
$ws.Range("N82").Value = -2503.6666

$ws.Range("H85").Value = 2233439.8
$ws.Range("J85").Value = 1781.6666
$ws.Range("L85").Value = 1781.6666
$ws.Range("N85").Value = -4277.6666

$ws.Range("H113").Value = 4951.4287
$ws.Range("I113").Value = 5226.6665
$ws.Range("K113").Value = 5226.6665
$ws.Range("M113").Value = -3056.6665

$ws.Range("H132").Value = 3737
$ws.Range("I132").Value = 2638.2
$ws.Range("K132").Value = 7914.599999999999
$ws.Range("M132").Value = -5384.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1669460.1
$ws.Range("I100").Value = 5001740.5
$ws.Range("K100").Value = 10003481
$ws.Range("M100").Value = -10002940

$ws.Range("H107").Value = 1140.6522
$ws.Range("I107").Value = 1250.875
$ws.Range("J107").Value = 888.7143
$ws.Range("K107").Value = 3752.625
$ws.Range("L107").Value = 2666.1429
$ws.Range("M107").Value = -1832.625
$ws.Range("N107").Value = -6506.1429

$ws.Range("H126").Value = 1942.091
$ws.Range("I126").Value = 1870.75
$ws.Range("J126").Value = 1982.8572
$ws.Range("K126").Value = 5612.25
$ws.Range("L126").Value = 5948.571599999999
$ws.Range("M126").Value = -3142.25
$ws.Range("N126").Value = -10888.5716

$ws.Range("H136").Value = 7447.633
$ws.Range("I136").Value = 3113.3462
$ws.Range("J136").Value = 9012.791999999999
$ws.Range("K136").Value = 9340.0386
$ws.Range("L136").Value = 27038.376
$ws.Range("M136").Value = -6790.0386
$ws.Range("N136").Value = -32138.376
